$wb = $excel.ActiveWorkbook

# The "Repayment schedule" sheet becomes the active sheet (was "Transactions").
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

# Insert a new (blank) column before column N to make room for a "Variable
# Instalments" related value, shifting the former N/O/P columns to O/P/Q.
$ws.Columns("N").Insert()

# Match the width used by the neighbouring "In Advance" column (M).
$ws.Range("N1").ColumnWidth = 9.83

# Leave the selection where the author left it after the edit.
$null = $ws.Range("S5").Select()
